# Auto-generated edit script: updates market-data derived value columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# 8 job sheets, per a scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 49.642857
$ws.Range("I2").Value = 20
$ws.Range("K2").Value = 20
$ws.Range("M2").Value = 93
$ws.Range("H17").Value = 1319.3103
$ws.Range("H18").Value = 71447390
$ws.Range("I18").Value = 100023920
$ws.Range("K18").Value = 100023920
$ws.Range("M18").Value = -100023636
$ws.Range("H40").Value = 3910.9583
$ws.Range("J40").Value = 4723.25
$ws.Range("L40").Value = 4723.25
$ws.Range("N40").Value = -5073.25
$ws.Range("H88").Value = 20872574
$ws.Range("J88").Value = 60723.9
$ws.Range("L88").Value = 60723.9
$ws.Range("N88").Value = -61535.9
$ws.Range("H91").Value = 20872574
$ws.Range("J91").Value = 60723.9
$ws.Range("L91").Value = 60723.9
$ws.Range("N91").Value = -63531.9
$ws.Range("H93").Value = 49268.8
$ws.Range("J93").Value = 49268.8
$ws.Range("L93").Value = 49268.8
$ws.Range("N93").Value = -54260.8
$ws.Range("H100").Value = 2775.3333
$ws.Range("I100").Value = 1502.25
$ws.Range("K100").Value = 1502.25
$ws.Range("M100").Value = -961.25
$ws.Range("H107").Value = 57503870
$ws.Range("I107").Value = 45002340
$ws.Range("K107").Value = 45002340
$ws.Range("M107").Value = -45000420
$ws.Range("H116").Value = 25007394
$ws.Range("I116").Value = 125000980
$ws.Range("J116").Value = 8999.875
$ws.Range("K116").Value = 125000980
$ws.Range("L116").Value = 8999.875
$ws.Range("M116").Value = -124997538
$ws.Range("N116").Value = -15883.875
$ws.Range("H137").Value = 3006.7222
$ws.Range("J137").Value = 2838.6667
$ws.Range("L137").Value = 8516.000100000001
$ws.Range("N137").Value = -13616.0001
$ws.Range("H138").Value = 7200.1953
$ws.Range("I138").Value = 2023.875
$ws.Range("J138").Value = 8455.061
$ws.Range("K138").Value = 6071.625
$ws.Range("L138").Value = 25365.183
$ws.Range("M138").Value = -931.625
$ws.Range("N138").Value = -35645.183

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 47621650
$ws.Range("I2").Value = 891.0769
$ws.Range("K2").Value = 891.0769
$ws.Range("M2").Value = -778.0769
$ws.Range("H19").Value = 200
$ws.Range("I19").Value = 300
$ws.Range("K19").Value = 300
$ws.Range("M19").Value = -71
$ws.Range("H32").Value = 1673596.9
$ws.Range("I32").Value = 1766264.1
$ws.Range("K32").Value = 1766264.1
$ws.Range("M32").Value = -1765977.1
$ws.Range("H57").Value = 4944.091
$ws.Range("I57").Value = 4944.091
$ws.Range("K57").Value = 4944.091
$ws.Range("M57").Value = -4460.091
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()  # was -52954
$ws.Range("H116").Value = 47621650
$ws.Range("I116").Value = 891.0769
$ws.Range("K116").Value = 891.0769
$ws.Range("M116").Value = 1402.9231
$ws.Range("H122").Value = 2874.7742
$ws.Range("I122").Value = 1581.5
$ws.Range("K122").Value = 4744.5
$ws.Range("M122").Value = -2294.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 47621650
$ws.Range("I3").Value = 891.0769
$ws.Range("K3").Value = 891.0769
$ws.Range("M3").Value = -777.0769
$ws.Range("H20").Value = 6174630
$ws.Range("I20").Value = 9261062
$ws.Range("K20").Value = 9261062
$ws.Range("M20").Value = -9260815
$ws.Range("H22").Value = 299
$ws.Range("I22").Value = 299
$ws.Range("K22").Value = 299
$ws.Range("M22").Value = -126
$ws.Range("H25").Value = 1900
$ws.Range("I25").Value = 1900
$ws.Range("K25").Value = 1900
$ws.Range("M25").Value = -1665
$ws.Range("H107").Value = 45002864
$ws.Range("I107").Value = 59211720
$ws.Range("K107").Value = 59211720
$ws.Range("M107").Value = -59209800

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5525.026
$ws.Range("J31").Value = 10951.556
$ws.Range("L31").Value = 10951.556
$ws.Range("N31").Value = -11541.556
$ws.Range("H33").Value = 1950.5
$ws.Range("I33").Value = 1950.5
$ws.Range("K33").Value = 1950.5
$ws.Range("M33").Value = -1571.5
$ws.Range("H34").Value = 5525.026
$ws.Range("J34").Value = 10951.556
$ws.Range("L34").Value = 10951.556
$ws.Range("N34").Value = -11355.556
$ws.Range("H50").Value = 45000
$ws.Range("J50").Value = 45000
$ws.Range("L50").Value = 45000
$ws.Range("N50").Value = -46250
$ws.Range("H58").Value = 15632968
$ws.Range("I58").Value = 38464420
$ws.Range("J58").Value = 11447.211
$ws.Range("K58").Value = 38464420
$ws.Range("L58").Value = 11447.211
$ws.Range("M58").Value = -38464217
$ws.Range("N58").Value = -11853.211
$ws.Range("H76").Value = 5006.857
$ws.Range("I76").Value = 5006.857
$ws.Range("K76").Value = 5006.857
$ws.Range("M76").Value = -4691.857
$ws.Range("H79").Value = 5006.857
$ws.Range("I79").Value = 5006.857
$ws.Range("K79").Value = 5006.857
$ws.Range("M79").Value = -3914.857
$ws.Range("H99").Value = 8499
$ws.Range("I99").Value = 6996
$ws.Range("J99").Value = 9000
$ws.Range("K99").Value = 6996
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = -5498
$ws.Range("N99").Value = -11996
$ws.Range("H126").Value = 8499
$ws.Range("I126").Value = 6996
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 20988
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -18518
$ws.Range("N126").Value = -31940
$ws.Range("H132").Value = 5702.9775
$ws.Range("I132").Value = 3050.3076
$ws.Range("J132").Value = 9332.947
$ws.Range("K132").Value = 9150.9228
$ws.Range("L132").Value = 27998.841
$ws.Range("M132").Value = -6620.9228
$ws.Range("N132").Value = -33058.841
$ws.Range("H134").Value = 4119.466
$ws.Range("I134").Value = 1845.7391
$ws.Range("J134").Value = 7993.222
$ws.Range("K134").Value = 5537.2173
$ws.Range("L134").Value = 23979.666
$ws.Range("M134").Value = -3002.2173
$ws.Range("N134").Value = -29049.666
$ws.Range("H136").Value = 15632968
$ws.Range("I136").Value = 38464420
$ws.Range("J136").Value = 11447.211
$ws.Range("K136").Value = 115393260
$ws.Range("L136").Value = 34341.633
$ws.Range("M136").Value = -115390710
$ws.Range("N136").Value = -39441.633

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 111403670
$ws.Range("I9").Value = 22949.5
$ws.Range("J9").Value = 155955950
$ws.Range("K9").Value = 68848.5
$ws.Range("L9").Value = 467867850
$ws.Range("M9").Value = -68624.5
$ws.Range("N9").Value = -467868298
$ws.Range("H17").Value = 541.375
$ws.Range("I17").Value = 81
$ws.Range("J17").Value = 607.1429000000001
$ws.Range("K17").Value = 243
$ws.Range("L17").Value = 1821.4287
$ws.Range("M17").Value = -74
$ws.Range("N17").Value = -2159.4287

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2019361.5
$ws.Range("I122").Value = 3159220.5
$ws.Range("J122").Value = 2687.7693
$ws.Range("K122").Value = 9477661.5
$ws.Range("L122").Value = 8063.3079
$ws.Range("M122").Value = -9475211.5
$ws.Range("N122").Value = -12963.3079

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 6984.5
$ws.Range("I43").Value = 6984.5
$ws.Range("K43").Value = 6984.5
$ws.Range("M43").Value = -6791.5
$ws.Range("H47").Value = 18840.5
$ws.Range("J47").Value = 18840.5
$ws.Range("L47").Value = 18840.5
$ws.Range("N47").Value = -19820.5
$ws.Range("H52").Value = 18840.5
$ws.Range("J52").Value = 18840.5
$ws.Range("L52").Value = 18840.5
$ws.Range("N52").Value = -19306.5
$ws.Range("H56").Value = 29333.334
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()  # was -26182
$ws.Range("H58").Value = 2789.5
$ws.Range("I58").Value = 2629
$ws.Range("K58").Value = 2629
$ws.Range("M58").Value = -2369
$ws.Range("H61").Value = 6579.5625
$ws.Range("I61").Value = 2799.6
$ws.Range("K61").Value = 2799.6
$ws.Range("M61").Value = -2597.6
$ws.Range("H113").Value = 6579.5625
$ws.Range("I113").Value = 2799.6
$ws.Range("K113").Value = 2799.6
$ws.Range("M113").Value = -629.5999999999999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 64769.5
$ws.Range("J92").Value = 64769.5
$ws.Range("L92").Value = 64769.5
$ws.Range("N92").Value = -69761.5
$ws.Range("H122").Value = 194800.67
$ws.Range("I122").Value = 365854.38
$ws.Range("J122").Value = 6641.6
$ws.Range("K122").Value = 1097563.14
$ws.Range("L122").Value = 19924.8
$ws.Range("M122").Value = -1095113.14
$ws.Range("N122").Value = -24824.8
$ws.Range("H132").Value = 21758068
$ws.Range("I132").Value = 45466636
$ws.Range("J132").Value = 25212.416
$ws.Range("K132").Value = 136399908
$ws.Range("L132").Value = 75637.24800000001
$ws.Range("M132").Value = -136397378
$ws.Range("N132").Value = -80697.24800000001

